$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''61.957.98'
$ws.Range("E2").Value = '''  -2.27%  '
$ws.Range("D3").Value = '''2.497.03'
$ws.Range("E3").Value = '''  -3.57%  '
$ws.Range("E4").Value = '''  +0.03%  '
$ws.Range("D5").Value = '''550.16'
$ws.Range("E5").Value = '''  -3.80%  '
$ws.Range("D6").Value = '''147.08'
$ws.Range("E6").Value = '''  -5.22%  '
$ws.Range("E7").Value = '''  +0.02%  '
$ws.Range("D8").Value = '''0.599'
$ws.Range("E8").Value = '''  -3.60%  '
$ws.Range("D9").Value = '''2.493.80'
$ws.Range("E9").Value = '''  -3.58%  '
$ws.Range("E10").Value = '''  -9.28%  '
$ws.Range("E11").Value = '''  -1.65%  '
$ws.Range("D12").Value = '''5.38'
$ws.Range("E12").Value = '''  -8.00%  '
$ws.Range("D13").Value = '''0.355'
$ws.Range("E13").Value = '''  -6.60%  '
$ws.Range("D14").Value = '''26.09'
$ws.Range("E14").Value = '''  -7.36%  '
$ws.Range("D15").Value = '''2.951.25'
$ws.Range("E15").Value = '''  -3.43%  '
$ws.Range("D16").Value = '''61.861.65'
$ws.Range("E16").Value = '''  -2.18%  '
$ws.Range("D17").Value = '''0.0000163'
$ws.Range("E17").Value = '''  -8.55%  '
$ws.Range("D18").Value = '''2.504.21'
$ws.Range("E18").Value = '''  -2.91%  '
$ws.Range("D19").Value = '''11.13'
$ws.Range("E19").Value = '''  -7.04%  '
$ws.Range("D21").Value = '''4.18'
$ws.Range("E21").Value = '''  -7.69%  '
$ws.Range("D22").Value = '''322.13'
$ws.Range("E22").Value = '''  -5.88%  '
$ws.Range("D23").Value = '''0.999'
$ws.Range("E23").Value = '''  -0.05%  '
$ws.Range("D24").Value = '''63.96'
$ws.Range("E24").Value = '''  -4.81%  '
$ws.Range("E25").Value = '''  -4.27%  '
$ws.Range("D26").Value = '''0.0000102'
$ws.Range("E26").Value = '''  -5.66%  '
$ws.Range("D27").Value = '''2.628.42'
$ws.Range("E27").Value = '''  -3.12%  '
$ws.Range("D28").Value = '''542.07'
$ws.Range("E28").Value = '''  -6.26%  '
$ws.Range("D29").Value = '''1.49'
$ws.Range("E29").Value = '''  -4.34%  '
$ws.Range("E30").Value = '''  -0.06%  '
$ws.Range("D31").Value = '''8.35'
$ws.Range("E31").Value = '''  -8.51%  '
$ws.Range("D32").Value = '''7.68'
$ws.Range("E32").Value = '''  -2.58%  '
$ws.Range("D33").Value = '''0.148'
$ws.Range("E33").Value = '''  -7.81%  '
$ws.Range("D34").Value = '''1.89'
$ws.Range("E34").Value = '''  -7.84%  '
$ws.Range("E35").Value = '''  -8.42%  '
$ws.Range("D36").Value = '''5.90'
$ws.Range("E36").Value = '''  -10.05%  '
$ws.Range("D37").Value = '''4.86'
$ws.Range("E37").Value = '''  -10.81%  '
$ws.Range("D38").Value = '''0.999'
$ws.Range("E38").Value = '''  +0.09%  '
$ws.Range("D39").Value = '''0.379'
$ws.Range("E39").Value = '''  -5.76%  '
$ws.Range("D40").Value = '''18.53'
$ws.Range("E40").Value = '''  -5.94%  '
$ws.Range("D41").Value = '''142.90'
$ws.Range("E41").Value = '''  -7.56%  '
$ws.Range("D42").Value = '''1.00'
$ws.Range("E42").Value = '''  +0.04%  '
$ws.Range("E43").Value = '''  -8.64%  '
$ws.Range("D44").Value = '''40.51'
$ws.Range("E44").Value = '''  -1.82%  '
$ws.Range("D45").Value = '''2.32'
$ws.Range("E45").Value = '''  -5.99%  '
$ws.Range("D46").Value = '''149.35'
$ws.Range("E46").Value = '''  -4.15%  '
$ws.Range("D47").Value = '''3.57'
$ws.Range("E47").Value = '''  -8.84%  '
$ws.Range("D48").Value = '''20.88'
$ws.Range("E48").Value = '''  -9.54%  '
$ws.Range("D49").Value = '''0.0534'
$ws.Range("E49").Value = '''  -9.07%  '
$ws.Range("D50").Value = '''0.590'
$ws.Range("E50").Value = '''  -5.34%  '
$ws.Range("D51").Value = '''0.0942'
$ws.Range("E51").Value = '''  -6.10%  '
